# Paginas verificadas y alerta incluida
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Add an (empty, default-styled) cell in column C for every existing
#    data row (1-24). These become plain "Normal" styled blank cells.
$ws.Range("C1:C24").Style = "Normal"

# 2) Row 25 gets a new "Combo box Del country" alternate value in column C:
#    "UKRAINE" (uses the same look as the rest of column A/worksheet text).
$ws.Range("C25").Value = "UKRAINE"

# 3) New row 26: a single underlined note cell in column B warning that this
#    particular value does not trigger any alert.
$ws.Range("B26").Value = "Con este valor no se abre ninguna alerta"
$ws.Range("B26").Font.Name = "Ubuntu"
$ws.Range("B26").Font.Size = 13
$ws.Range("B26").Font.Color = 0
$ws.Range("B26").Font.Underline = $true

# Update the active selection to match the row that was just edited.
$ws.Range("B26").Select()

Write-Host "edit applied"
